# Merge the two runs "HTML, CSS, JavaScript, " + "GitHub" into a single
# run "HTML, CSS, JavaScript, GitHub" inside the "DOCTOR BOOKING WEBSITE"
# project textbox on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the textbox shape holding the project description
# ("TextBox 88" in the authored deck).
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 88") {
        $shape = $candidate
        break
    }
}

$textFrame = $shape.TextFrame
$textRange = $textFrame.TextRange

# The second paragraph currently holds the tech-stack line, split across
# two runs ("HTML, CSS, JavaScript, " and "GitHub") with identical
# formatting. Remember the shape's current (auto-fitted) height so we can
# restore it afterwards, since re-assigning paragraph text causes the
# auto-fit textbox to recompute its height as a side effect.
$origHeight = $shape.Height

$paragraph = $textRange.Paragraphs(2, 1)

# Re-assigning .Text with content that shares a common prefix/suffix with
# the existing runs only patches the differing runs in place (preserving
# the two-run split). Routing through an unrelated placeholder value first
# forces the host to collapse the paragraph into a single run (taking on
# the formatting of the first original run), after which we set the final
# desired text - yielding one run "HTML, CSS, JavaScript, GitHub".
$paragraph.Text = "placeholder-for-run-merge"
$paragraph.Text = "HTML, CSS, JavaScript, GitHub"

# Undo the auto-fit height side effect triggered by the text edits above.
$shape.Height = $origHeight
